# Horarios actualizados Línea 141 - 1056
# Refresh the scraped "last updated" timestamp, row counts and minute
# countdowns across the three worksheets, and append the new arrival
# (11_ETCHEVERRY) to the main "LP1912" sheet.

$wb = $excel.ActiveWorkbook

$oldTs = "02:48:47"
$newTs = "02:59:45"

# ---------------------------------------------------------------------
# Sheet 1: "LP1912" - full detail sheet, gains a new row (10)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTs"
$ws1.Range("A3").Value = "Total filas: 5"

$ws1.Cells.Item(6, 1).Value = $newTs
$ws1.Cells.Item(6, 2).Value = "02:59"
$ws1.Cells.Item(6, 4).Value = 0

$ws1.Cells.Item(7, 1).Value = $newTs
$ws1.Cells.Item(7, 4).Value = 49

$ws1.Cells.Item(8, 1).Value = $newTs
$ws1.Cells.Item(8, 4).Value = 62

$ws1.Cells.Item(9, 1).Value = $newTs
$ws1.Cells.Item(9, 4).Value = 106

# New row 10: next scheduled arrival
$ws1.Cells.Item(10, 1).Value = $newTs
$ws1.Cells.Item(10, 2).Value = "04:53"
$ws1.Cells.Item(10, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(10, 4).Value = 114
$ws1.Cells.Item(10, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: "LP1912-215" - filtered detail sheet (no new rows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTs"

$ws2.Cells.Item(6, 1).Value = $newTs
$ws2.Cells.Item(6, 2).Value = "02:59"
$ws2.Cells.Item(6, 4).Value = 0

$ws2.Cells.Item(7, 1).Value = $newTs
$ws2.Cells.Item(7, 4).Value = 106

# ---------------------------------------------------------------------
# Sheet 3: "6203-6173" - summary-only sheet (just the timestamp)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTs"
